# Update the "doctor_MA" / average column (AF) values on Sheet1
# per updated results described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AF4").Value  = 0.755
$ws.Range("AF5").Value  = 0.971
$ws.Range("AF6").Value  = 0.849
$ws.Range("AF7").Value  = 0.918
$ws.Range("AF8").Value  = 0.88
$ws.Range("AF9").Value  = 0.735
$ws.Range("AF10").Value = 0.971
$ws.Range("AF11").Value = 0.971
$ws.Range("AF12").Value = 1.273
$ws.Range("AF13").Value = 1.559
